# The worksheet currently has an empty row 2 (data starts at row 1 header,
# then jumps to row 3). This edit removes that blank row so the data block
# becomes contiguous (rows 1-5 instead of 1,3-6), and fixes the G column
# value for the "ATIF ALABBASI" record so it stores the plain text name
# instead of a stray numeric value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the empty second row; this shifts the old rows 3,4,5,6 up to 2,3,4,5.
$ws.Rows.Item(2).Delete()

# Replace the numeric placeholder in G4 (old G5, row for ATIF ALABBASI) with
# the plain text value so it just holds the name, matching column A.
$ws.Range("G4").Value = "ATIF ALABBASI"

# Update the worksheet dimension/selection to match the new, smaller data range.
$ws.Range("A1:G5").Select()
